# PM11 Tidsregistrering for Rasmus.xlsx - update time entries
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: end time (E22) corrected from 11:30 to 11:00 -> cascades through
# the running-total formulas in G/H for the rest of the sheet.
$ws.Range("E22").Value = 0.45833333333333331

# Row 25: F25 ("1 timer") cell alignment changed (now centered, no longer
# vertically centered) - matches the style used on the new F26 entry below.
$ws.Range("F25").HorizontalAlignment = -4108

# Row 26: new time entry - "implementering af SD0104 med Sune"
$ws.Range("A26").Value = "implementering af SD0104 med Sune"
$ws.Range("B26").Value = "Implenter"
$ws.Range("C26").Value = 43892
$ws.Range("D26").Value = 0.58333333333333337
$ws.Range("E26").Value = 0.625
$ws.Range("F26").Value = "1 timer"
$ws.Range("F26").HorizontalAlignment = -4108

# Row 27: new time entry - "Review af UC07"
$ws.Range("A27").Value = "Review af UC07"
$ws.Range("B27").Value = "Reviewer"
$ws.Range("C27").Value = 43893
$ws.Range("D27").Value = 0.375
$ws.Range("E27").Value = 0.39583333333333331
$ws.Range("F27").Value = "30 minutter"

# View state: zoom bumped slightly and selection moved to the newly
# entered row.
$ws.Range("F27").Select()
$excel.ActiveWindow.Zoom = 68
